$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C ("Spare part name"), which
# pushes that column and everything to its right two columns over.
$ws.Range("C:D").Insert() | Out-Null

# Give the two newly inserted columns their header text in row 2. The
# insert already carried column B's header formatting (bold text on a
# yellow fill, centered) onto the new C/D cells.
$ws.Range("C2").Value = "Internal Use"
$ws.Range("D2").Value = "Warehouse"

# Match the column widths recorded for the new layout as closely as
# this runtime's column-width rounding allows.
$ws.Range("C:C").ColumnWidth = 15.666666666666666
$ws.Range("D:D").ColumnWidth = 16.333333333333332

# Rebuild the autofilter so it spans the table's new extent (A2:Q2)
# instead of the stale A2:O2 range left over from the insert.
$ws.AutoFilterMode = $false
$ws.Range("A2:Q2").AutoFilter() | Out-Null

# Keep the workbook's hidden _FilterDatabase name in sync with the
# filter's new range.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$2:`$Q`$2"
